# Select the "Diameter Premiums" worksheet and insert a new column A
# ("SHAPE") ahead of the existing data, filling it with "RO" for every
# data row. All previously existing columns (CUT, POLY, SYM, FLUO, SIZE,
# DIAMETER_MIN, DIAMETER_MAX, KEY_COLOR_CLARITY, DISCOUNT) shift one
# column to the right (A->B ... I->J).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diameter Premiums")

# Insert a new, blank column before column A; this shifts existing data
# (and formatting) from columns A:I to B:J.
$ws.Columns.Item(1).Insert()

# The new header cell (A1) needs the same formatting as the header row it
# now belongs to (bordered, bold, centered) -- copy that format from the
# neighboring header cell B1 (the old A1 / "CUT" header).
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Find the last used row in the sheet (should be 85).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Header for the newly inserted column.
$ws.Cells.Item(1, 1).Value = "SHAPE"

# Fill the new column's data rows with "RO".
$dataRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 1))
$dataRange.Value = "RO"
